$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6226591760299626
$ws1.Range("C2").Value = 0.5711183496199783
$ws1.Range("D2").Value = 0.9850187265917603
$ws1.Range("E2").Value = 0.7230240549828179
$ws1.Range("F2").Value = 0.8603205757278377
$ws1.Range("G2").Value = 0.958307056267956
$ws1.Range("H2").Value = 0.7297689685645752
$ws1.Range("I2").Value = 526
$ws1.Range("J2").Value = 395
$ws1.Range("K2").Value = 139
$ws1.Range("L2").Value = 8

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9455782312925171
$ws2.Range("C2").Value = 0.2602996254681648
$ws2.Range("D2").Value = 0.408223201174743

$ws2.Range("B3").Value = 0.5711183496199783
$ws2.Range("C3").Value = 0.9850187265917603
$ws2.Range("D3").Value = 0.7230240549828179

$ws2.Range("B4").Value = 0.6226591760299626
$ws2.Range("C4").Value = 0.6226591760299626
$ws2.Range("D4").Value = 0.6226591760299626
$ws2.Range("E4").Value = 0.6226591760299626

$ws2.Range("B5").Value = 0.7583482904562477
$ws2.Range("C5").Value = 0.6226591760299626
$ws2.Range("D5").Value = 0.5656236280787804

$ws2.Range("B6").Value = 0.7583482904562476
$ws2.Range("C6").Value = 0.6226591760299626
$ws2.Range("D6").Value = 0.5656236280787804

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 139
$ws3.Range("C2").Value = 395
$ws3.Range("B3").Value = 8
$ws3.Range("C3").Value = 526
